$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 903, shifting the existing row 903 (and everything
# below it) down by one. This grows the used range from A1:D944 to A1:D945.
$ws.Rows.Item(903).Insert()

# Populate the newly-inserted row 903 with the new data point for
# 2026/02/28 (continuing the existing 土 values for that date).
# Column A holds dates formatted as plain text, so force text formatting
# before assigning the value to stop Excel auto-converting the
# date-looking string into a real date serial number, then clear the
# formatting again so the cell matches its neighbours (no explicit style).
$ws.Range("A903").NumberFormat = "@"
$ws.Range("A903").Value = "2026/02/28"
$ws.Range("A903").ClearFormats()

$ws.Range("B903").Value = "土"
$ws.Range("C903").Value = 16
$ws.Range("D903").Value = 201
